$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update the two input values (Gross Expenditures and Total Labor Cost)
$ws.Range("D3").Value = 279239.07
$ws.Range("D5").Value = 39314.870000000003

# Update the active selection to match the saved view state
$ws.Range("G7:I16").Select()
